$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 2-10 (cols B:G) down into rows 3-11,
# preserving the current values (a new quarter's data is being
# inserted at the top of the history, pushing the older rows down).
for ($r = 10; $r -ge 2; $r--) {
    $src = $ws.Range("B" + $r + ":G" + $r)
    $dst = $ws.Range("B" + ($r + 1) + ":G" + ($r + 1))
    $dst.Value = $src.Value()
}

# Populate row 2 with the newly computed values for this quarter.
$ws.Range("B2").Value = 0.1783908196033299
$ws.Range("C2").Value = 0.3606156554386025
$ws.Range("D2").Value = 0.2599511937740667
$ws.Range("E2").Value = 0.5098540906711122
$ws.Range("F2").Value = 0.4943913024279584
$ws.Range("G2").Value = 15
